$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sites")
$ws.Columns.Item(2).Insert()
$ws.Range("B2:B28").Formula = "=VLOOKUP(A2, 'Site listes'!A2:C28, 3, FALSE)"
$ws.Range("B1").Value = "SiteNew"
$ws.Range("C1").Value = "GPS"
$ws.Range("B1:B28").HorizontalAlignment = -4108
$ws.Range("B1").Borders.Item(7).LineStyle = -4142
$ws.Range("B2:B28").Borders.Item(7).LineStyle = -4142
Write-Output "done"
